# Apply the "Added a few more slots" edit to the Calaveras Explosivas
# review document:
#   1. Insert a new "Meta description: ..." paragraph right after the
#      Heading1 title paragraph.
#   2. Remove the stray duplicate bold "Play Calaveras Explosivas for
#      free - gameplay, symbols, and RTP" paragraph near the end of the
#      document.
#   3. Replace the text of the remaining (italic) paragraph at the very
#      end with the new feature-image generation prompt, keeping its
#      italic formatting intact.

$d = $word.ActiveDocument

# --- 1. Insert the "Meta description" paragraph after the title -----------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the unique exploding symbols gameplay and themed design of Calaveras Explosivas slot game for free. Discover its low volatility and RTP value.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# --- 2. Delete the duplicated bold title paragraph near the end -----------
# Locate it via the paragraph collection: the duplicate is the
# second-to-last paragraph in the document, holding only the bold title.
# Search from the tail backwards (stopping before paragraph 1, the real
# Heading1 title) so the first hit is always that trailing duplicate.
$count = $d.Paragraphs.Count
$dupPara = $null
for ($i = $count; $i -ge 2; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "Play Calaveras Explosivas for free - gameplay, symbols, and RTP") {
        $dupPara = $para
        break
    }
}
if ($dupPara -ne $null) {
    $dupPara.Range.Delete()
}

# --- 3. Replace the final italic paragraph's text --------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Find.Execute(
    "Explore the unique exploding symbols gameplay and themed design of Calaveras Explosivas slot game for free. Discover its low volatility and RTP value.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a feature image for Calaveras Explosivas: Design a colorful cartoon-style image of a happy Maya warrior wearing glasses. The warrior should have a big smile on his face and be surrounded by exploding skulls, the golden wild symbol, and the white skull with roses in place of the eyes scatter symbol. The background should feature the typical street in a Mexican village with flags running across from one balcony to another. The image should be eye-catching and reflect the festive and lively atmosphere of the Calaveras Explosivas slot game.",
    2)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
